# fixed a bug in rtp
# Rewrites the data rows (2-21) of the active worksheet with their corrected
# values/order. Row 1 (headers) and rows 22-26 (totals area) are untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(601, 9, 60, 67, 60, 42),
    @(201, 9, 30, 15, 45, 30),
    @(1203, 3, 15, 15, 15, 15),
    @(901, 16, 15, 45, 60, 60),
    @(301, 6, 45, 30, 60, 45),
    @(401, 9, 48, 67, 75, 45),
    @(801, 3, 67, 65, 52, 45),
    @(1202, 2, 10, 10, 10, 10),
    @(902, 1, 0, 0, 0, 0),
    @(501, 9, 52, 30, 75, 45),
    @(701, 3, 90, 45, 97, 15),
    @(101, 9, 30, 15, 60, 15),
    @(1201, 2, 10, 10, 10, 10),
    @(1001, 18, 30, 75, 60, 72),
    @(502, 0, 4, 0, 0, 0),
    @(1, 0, 2, 2, 2, 2),
    @(1101, 0, 15, 30, 30, 0),
    @(3, 0, 3, 3, 3, 3),
    @(802, 0, 4, 5, 4, 0),
    @(2, 0, 2, 2, 2, 2)
)

$startRow = 2
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $values = $data[$i]
    for ($col = 1; $col -le $values.Count; $col++) {
        $ws.Cells.Item($row, $col).Value = $values[$col - 1]
    }
}
